# This script reproduces, on the "展览" (Exhibitions) and "全部类型" (All
# types) sheets, the same logical edit: the "苏州·代号鸢only茶话会-星渡咖啡"
# listing is removed from its original slot (2024-08-03) and re-added later
# in the list with updated dates/price near the very end (2024-10-19,
# right before the last "华盟国漫次元嘉年华" row); a brand-new listing
# "苏州·艾卡动漫游戏嘉年华（免票展）" is inserted before "太仓·第六届龙狮动漫嘉年华";
# and a handful of "want-to-go" counts (column F) on unrelated rows are
# bumped up.

$wb = $excel.ActiveWorkbook

function Update-ExpoSheet(
    $ws,
    $rowOsuju,        # row holding 苏州·授渔仲夏动漫节2.0 (after the delete-shift settles it)
    $rowOcgCarnival,  # row holding 苏州·第三届.OCG.Summer Carnival-国潮动漫游戏嘉年华
    $rowTaicang,      # row holding 太仓·第六届龙狮动漫嘉年华 (insertion point for 艾卡)
    $rowIcan,         # row holding 苏州·ICAN summer World动漫品牌夏游节
    $rowRedamancy,    # row holding 苏州·第二届Redamancy动漫游戏嘉年华
    $rowVolleyball,   # row holding 苏州·排球少年only-茶歇
    $rowUmamusume,    # row holding 苏州·赛马娘ONLY
    $rowGoodJump,     # row holding 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
    $rowIcome,        # row holding 苏州·I COME ACG动漫品牌博览会
    $rowLixiangxiang, # row holding 苏州·第十三届理想乡动漫展-同人创作者大会
    $lixiangxiangNewValue,
    $rowOcg4th,       # row holding 苏州·第四届-OCG国朝动漫游戏嘉年华
    $rowBaihe         # row holding 苏州·第二届百合Only同人展交流
) {
    # --- simple same-row numeric bumps (rows 3 and 4 stay put) ---
    $ws.Range("F3").Value = 511
    $ws.Range("F4").Value = 1501

    # --- remove the "代号鸢only茶话会" row; everything below shifts up ---
    $ws.Rows.Item(7).Delete()

    # --- bumps that land on rows which shifted up by one ---
    $ws.Range("F" + $rowOsuju).Value = 729
    $ws.Range("F" + $rowOcgCarnival).Value = 6343

    # --- insert the brand-new "艾卡动漫游戏嘉年华（免票展）" row ---
    $ws.Rows.Item($rowTaicang).Insert()
    $rowAika = $rowTaicang
    $ws.Range("A" + ($rowAika - 1)).Copy()
    $ws.Range("A" + $rowAika).PasteSpecial(-4122)
    $ws.Range("A" + $rowAika).Value = ($rowAika - 1)
    $ws.Range("B" + $rowAika).Value = "2024-08-10"
    $ws.Range("C" + $rowAika).Value = "苏州·艾卡动漫游戏嘉年华（免票展）"
    $ws.Range("D" + $rowAika).Value = "相城大道3188号 苏州高铁吾悦广场"
    $ws.Range("E" + $rowAika).Value = "2024.08.10 13:00-08.10 18:00"
    $ws.Range("F" + $rowAika).Value = 0
    $ws.Range("G" + $rowAika).Value = 49
    $ws.Range("H" + $rowAika).Value = "https://show.bilibili.com/platform/detail.html?id=90010"
    $ws.Range("I" + $rowAika).Value = "//i0.hdslb.com/bfs/openplatform/202407/z3fQNCvo1722252927908.jpeg"

    # --- remaining numeric "want-to-go" bumps, unaffected by row identity ---
    $ws.Range("F" + $rowIcan).Value = 15204
    $ws.Range("F" + $rowRedamancy).Value = 1507
    $ws.Range("F" + $rowVolleyball).Value = 274
    $ws.Range("F" + $rowUmamusume).Value = 99
    $ws.Range("F" + $rowGoodJump).Value = 11007
    $ws.Range("F" + $rowIcome).Value = 739
    $ws.Range("F" + $rowLixiangxiang).Value = $lixiangxiangNewValue
    $ws.Range("F" + $rowOcg4th).Value = 226
    $ws.Range("F" + $rowBaihe).Value = 14

    # --- append the updated "代号鸢only茶话会" row right before the final row ---
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $ws.Rows.Item($lastRow).Insert()
    $rowDaihaoyuan = $lastRow
    $ws.Range("A" + ($rowDaihaoyuan - 1)).Copy()
    $ws.Range("A" + $rowDaihaoyuan).PasteSpecial(-4122)
    $ws.Range("A" + $rowDaihaoyuan).Value = ($rowDaihaoyuan - 1)
    $ws.Range("B" + $rowDaihaoyuan).Value = "2024-10-19"
    $ws.Range("C" + $rowDaihaoyuan).Value = "苏州·代号鸢only茶话会-星渡咖啡"
    $ws.Range("D" + $rowDaihaoyuan).Value = "德必姑苏WE国际文化艺术中心6-102室渔郎桥浜路16号 星渡咖啡"
    $ws.Range("E" + $rowDaihaoyuan).Value = "2024.10.19 10:00-10.20 19:00"
    $ws.Range("F" + $rowDaihaoyuan).Value = 299
    $ws.Range("G" + $rowDaihaoyuan).Value = 50
    $ws.Range("H" + $rowDaihaoyuan).Value = "https://show.bilibili.com/platform/detail.html?id=87685"
    $ws.Range("I" + $rowDaihaoyuan).Value = "//i1.hdslb.com/bfs/openplatform/202406/eyHRVQuv1718780132754.jpeg"
}

# ---- Sheet "展览" (31 data rows before edit -> 32 after) ----
$wsExpo = $wb.Worksheets.Item("展览")
Update-ExpoSheet $wsExpo 9 14 15 20 21 22 24 25 26 27 4301 28 30

# ---- Sheet "全部类型" (34 data rows before edit -> 35 after) ----
$wsAll = $wb.Worksheets.Item("全部类型")
Update-ExpoSheet $wsAll 10 17 18 23 24 25 27 28 29 30 4302 31 33
